# Update cryptos list values per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.456.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.380.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.92"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.72"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.12"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.65%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.109"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.989"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.41%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.742.48"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.54"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.381.11"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.303.00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.68"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +16.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.32"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.08%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.75"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.37"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.87"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.58"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.28"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.47"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0957"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.53"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "167.64"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.90"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.75"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.95"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +13.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.05"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.97"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.74"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.59"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.01"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.865.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +13.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.97"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.15"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.86"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.30"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.61%  "
